$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51
# D-column values that look like plain numbers get a leading apostrophe
# so Excel stores them as text (matching the source data format), exactly
# as they were already stored as text in the workbook before this edit.
$ws.Range("D2").Value = "63.049.18"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").Value = "2.626.23"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'603.95"
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("D6").Value = "'146.13"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -1.06%  "

$ws.Range("D9").Value = "2.624.70"
$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("D10").Value = "'0.107"
$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").Value = "'5.59"
$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").Value = "'0.362"
$ws.Range("E13").Value = "  +1.79%  "

$ws.Range("D14").Value = "'27.14"
$ws.Range("E14").Value = "  -1.93%  "

$ws.Range("D15").Value = "3.094.65"
$ws.Range("E15").Value = "  -1.49%  "

$ws.Range("D16").Value = "62.922.78"
$ws.Range("E16").Value = "  -1.33%  "

$ws.Range("E17").Value = "  -1.98%  "

$ws.Range("D18").Value = "2.617.13"
$ws.Range("E18").Value = "  -4.56%  "

$ws.Range("D19").Value = "'11.26"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("D20").Value = "'4.46"
$ws.Range("E20").Value = "  +2.09%  "

$ws.Range("D21").Value = "'339.15"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("D22").Value = "'6.85"
$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("E24").Value = "  -4.44%  "

$ws.Range("D25").Value = "'66.46"
$ws.Range("E25").Value = "  -2.60%  "

$ws.Range("E26").Value = "  -2.17%  "

$ws.Range("D27").Value = "'1.52"
$ws.Range("E27").Value = "  -5.52%  "

$ws.Range("D28").Value = "'8.62"
$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("D29").Value = "'0.162"
$ws.Range("E29").Value = "  -2.66%  "

$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").Value = "'534.38"
$ws.Range("E31").Value = "  -3.12%  "

$ws.Range("D32").Value = "'7.87"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").Value = "'2.02"
$ws.Range("E33").Value = "  +1.37%  "

$ws.Range("E34").Value = "  -1.66%  "

$ws.Range("E35").Value = "  -2.16%  "

$ws.Range("E36").Value = "  +12.25%  "

$ws.Range("D37").Value = "'168.84"
$ws.Range("E37").Value = "  -3.64%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("E40").Value = "  -0.88%  "

$ws.Range("E41").Value = "  +6.50%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").Value = "'167.72"
$ws.Range("E43").Value = "  -2.64%  "

$ws.Range("D44").Value = "'3.73"
$ws.Range("E44").Value = "  -0.64%  "

$ws.Range("D45").Value = "'22.18"
$ws.Range("E45").Value = "  +2.65%  "

$ws.Range("D46").Value = "'0.0563"
$ws.Range("E46").Value = "  +3.04%  "

$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").Value = "'0.0956"

$ws.Range("D50").Value = "'18.51"
$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("D51").Value = "'1.75"
$ws.Range("E51").Value = "  -0.41%  "
